$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 534, shifting rows 534:601 down to 535:601
$ws.Rows.Item(534).Insert()

# Populate the newly inserted row 534 with the new data
$ws.Cells.Item(534, 1).Value = 3
$ws.Cells.Item(534, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(534, 3).Value = "Coquimbo"
$ws.Cells.Item(534, 4).Value = 45124
$ws.Cells.Item(534, 5).Value = 5
$ws.Cells.Item(534, 6).Value = 100112043
$ws.Cells.Item(534, 7).Value = "Pepino ensalada"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 115
$ws.Cells.Item(534, 11).Value = 9000
$ws.Cells.Item(534, 12).Value = 9500
$ws.Cells.Item(534, 13).Value = 9217
$ws.Cells.Item(534, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(534, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(534, 16).Value = 154
$ws.Cells.Item(534, 17).Value = 60
$ws.Cells.Item(534, 18).Value = "Hortaliza"
